$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row with "Pieza 3D" expense entry
$ws.Range("A9").Value = "Pieza 3D"
$ws.Range("B9").Value = 25

# Move selection to follow the newly added row (as reflected in the diff)
$ws.Range("A10").Select()

$wb.Save()
